$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellB = $ws.Cells.Item($row, 2)

    $oldValue = [string]$cellB.Value2
    $newNumber = [int64]$oldValue + 10010
    $newValue = $newNumber.ToString()

    $cellA.Value = "http://localhost:80/12p?name=$newValue"

    # The new Test Kit Number text is purely numeric-looking, so a plain
    # .Value write would get auto-coerced to a number (losing the original
    # text type). Force text entry like Excel does for a "Text"-formatted
    # cell, then restore General formatting so no stray number format is
    # left behind on the cell.
    $cellB.NumberFormat = "@"
    $cellB.Value = $newValue
    $cellB.Style = "Normal"
}
